$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New/updated text blocks introduced by this edit
$abstractD2 = 'The IUPHAR/BPS Guide to PHARMACOLOGY (GtoPdb, www.
guidetopharmacology.
org) and its precursor IUPHAR-DB, have captured expert-curated interactions between targets and ligands from selected papers in pharmacology and drug discovery since 2003. This resource continues to be developed in conjunction with the International Union of Basic and Clinical Pharmacology (IUPHAR) and the British Pharmacological Society (BPS).
 As previously described, our unique model of content selection and quality control is based on 96 target-class subcommittees comprising 512 scientists collaborating with in-house curators.
 This update describes content expansion, new features and interoperability improvements introduced in the 10 releases since August 2015. Our relationship matrix now describes ∼9000 ligands, ∼15 000 binding constants, ∼6000 papers and ∼1700 human proteins.
 As an important addition, we also introduce our newly funded project for the Guide to IMMUNOPHARMACOLOGY (GtoImmuPdb, www.
guidetoimmunopharmacology.
org).
 This has been ‘forked’ from the well-established GtoPdb data model and expanded into new types of data related to the immune system and inflammatory processes.
 This includes new ligands, targets, pathways, cell types and diseases for which we are recruiting new IUPHAR expert committees.
 Designed as an immunopharmacological gateway, it also has an emphasis on potential therapeutic interventions.
'
$authorsE2 = '[Simon D%Harding%NULL%0, Joanna L%Sharman%NULL%2, Joanna L%Sharman%NULL%0, Elena%Faccenda%NULL%1, Chris%Southan%NULL%1, Adam J%Pawson%NULL%1, Sam%Ireland%NULL%1, Alasdair J G%Gray%NULL%1, Liam%Bruce%NULL%1, Stephen P H%Alexander%NULL%1, Stephen%Anderton%NULL%1, Clare%Bryant%NULL%1, Anthony P%Davenport%NULL%1, Christian%Doerig%NULL%1, Doriano%Fabbro%NULL%1, Francesca%Levi-Schaffer%NULL%1, Michael%Spedding%NULL%1, Jamie A%Davies%jamie.davies@ed.ac.uk%1, NULL%NULL%NULL%0]'
$authorsE5 = '[Carl%Llor%carles.llor@urv.cat%0, Ana%Moragas%amoragasm@meditex.es%1, Carolina%Bayona%cbayona.tarte.ics@gencat.cat%1, Rosa%Morros%rmorros@idiapjgol.org%1, Helena%Pera%hpera@idiapjgol.org%1, Josep M%Cots%23465jcy@comb.cat%1, Yvonne%Fernández%yfernandez.tarte.ics@gencat.cat%1, Marc%Miravitlles%marcm@separ.es%1, Albert%Boada%aboadav.bcn.ics@gencat.cat%1]'

# Row 2: IUPHAR/BPS Guide to PHARMACOLOGY citation - refreshed abstract/authors/publisher
$ws.Range("D2").Value = $abstractD2
$ws.Range("E2").Value = $authorsE2
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = "Oxford University Press"

# Rows 3, 4, 7, 8: citations reset to Unknown Title/Abstract placeholders
$ws.Range("C3").Value = "Unknown Title"
$ws.Range("D3").Value = "Unknown Abstract"
$ws.Range("E3").Value = "[]"
$ws.Range("F3").Value = "not found"
$ws.Range("G3").Value = "N/A"
# Leading apostrophe forces text so '1970-01-01' isn't reinterpreted as a date serial
$ws.Range("H3").Value = "'1970-01-01"
$ws.Range("I3").Value = ""

$ws.Range("C4").Value = "Unknown Title"
$ws.Range("D4").Value = "Unknown Abstract"
$ws.Range("E4").Value = "[]"
$ws.Range("F4").Value = "not found"
$ws.Range("G4").Value = "N/A"
# Leading apostrophe forces text so '1970-01-01' isn't reinterpreted as a date serial
$ws.Range("H4").Value = "'1970-01-01"
$ws.Range("I4").Value = ""

$ws.Range("C7").Value = "Unknown Title"
$ws.Range("D7").Value = "Unknown Abstract"
$ws.Range("E7").Value = "[]"
$ws.Range("F7").Value = "not found"
$ws.Range("G7").Value = "N/A"
# Leading apostrophe forces text so '1970-01-01' isn't reinterpreted as a date serial
$ws.Range("H7").Value = "'1970-01-01"
$ws.Range("I7").Value = ""

$ws.Range("C8").Value = "Unknown Title"
$ws.Range("D8").Value = "Unknown Abstract"
$ws.Range("E8").Value = "[]"
$ws.Range("F8").Value = "not found"
$ws.Range("G8").Value = "N/A"
# Leading apostrophe forces text so '1970-01-01' isn't reinterpreted as a date serial
$ws.Range("H8").Value = "'1970-01-01"
$ws.Range("I8").Value = ""

# Rows 6, 9: citations reset (title/authors/ID/format/date), clear misc data
$ws.Range("C6").Value = "Unknown Title"
$ws.Range("E6").Value = "[]"
$ws.Range("F6").Value = "not found"
$ws.Range("G6").Value = "N/A"
# Leading apostrophe forces text so '1970-01-01' isn't reinterpreted as a date serial
$ws.Range("H6").Value = "'1970-01-01"
$ws.Range("J6").Value = ""

$ws.Range("C9").Value = "Unknown Title"
$ws.Range("E9").Value = "[]"
$ws.Range("F9").Value = "not found"
$ws.Range("G9").Value = "N/A"
# Leading apostrophe forces text so '1970-01-01' isn't reinterpreted as a date serial
$ws.Range("H9").Value = "'1970-01-01"
$ws.Range("J9").Value = ""

# Row 5: BioMed Central citation - refreshed authors/publisher
$ws.Range("E5").Value = $authorsE5
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = "BioMed Central"
